# Auto-generated Excel COM-interop script
# Updates pricing/profit columns (H:N) across multiple item sheets
# to match the scheduled-runner data refresh described in the commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 5216205.5
$ws.Range("I88").Value = 715.375
$ws.Range("J88").Value = 6761536
$ws.Range("K88").Value = 715.375
$ws.Range("L88").Value = 6761536
$ws.Range("M88").Value = -309.375
$ws.Range("N88").Value = -6762348
$ws.Range("H91").Value = 5216205.5
$ws.Range("I91").Value = 715.375
$ws.Range("J91").Value = 6761536
$ws.Range("K91").Value = 715.375
$ws.Range("L91").Value = 6761536
$ws.Range("M91").Value = 688.625
$ws.Range("N91").Value = -6764344
$ws.Range("H92").Value = 717.3182
$ws.Range("I92").Value = 498.92307
$ws.Range("J92").Value = 1032.7778
$ws.Range("K92").Value = 498.92307
$ws.Range("L92").Value = 1032.7778
$ws.Range("M92").Value = 749.0769299999999
$ws.Range("N92").Value = -3528.7778
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H99").Value = 294.33334
$ws.Range("I99").Value = 301.0909
$ws.Range("K99").Value = 903.2727
$ws.Range("M99").Value = 594.7273
$ws.Range("H100").Value = 15879.571
$ws.Range("I100").Value = 21101.4
$ws.Range("J100").Value = 2825
$ws.Range("K100").Value = 21101.4
$ws.Range("L100").Value = 2825
$ws.Range("M100").Value = -20560.4
$ws.Range("N100").Value = -3907
$ws.Range("H101").Value = 11613758
$ws.Range("I101").Value = 333495.66
$ws.Range("J101").Value = 45454544
$ws.Range("K101").Value = 1000486.98
$ws.Range("L101").Value = 136363632
$ws.Range("M101").Value = -998864.98
$ws.Range("N101").Value = -136366876
$ws.Range("H103").Value = 1187.4375
$ws.Range("I103").Value = 2217.1667
$ws.Range("J103").Value = 569.6
$ws.Range("K103").Value = 6651.500100000001
$ws.Range("L103").Value = 1708.8
$ws.Range("M103").Value = -6065.500100000001
$ws.Range("N103").Value = -2880.8
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H129").Value = 55556700
$ws.Range("I129").Value = 142858370
$ws.Range("J129").Value = 1091
$ws.Range("K129").Value = 428575110
$ws.Range("L129").Value = 3273
$ws.Range("M129").Value = -428570110
$ws.Range("N129").Value = -13273
$ws.Range("H132").Value = 11369552
$ws.Range("I132").Value = 13520007
$ws.Range("J132").Value = 2858.8572
$ws.Range("K132").Value = 40560021
$ws.Range("L132").Value = 8576.571599999999
$ws.Range("M132").Value = -40557491
$ws.Range("N132").Value = -13636.5716
$ws.Range("H138").Value = 3372.575
$ws.Range("I138").Value = 1767.88
$ws.Range("J138").Value = 4101.982
$ws.Range("K138").Value = 5303.64
$ws.Range("L138").Value = 12305.946
$ws.Range("M138").Value = -163.6400000000003
$ws.Range("N138").Value = -22585.946
$ws.Range("H141").Value = 2695.8262
$ws.Range("I141").Value = 1678.1538
$ws.Range("J141").Value = 8365.714
$ws.Range("K141").Value = 5034.4614
$ws.Range("L141").Value = 25097.142
$ws.Range("M141").Value = 145.5385999999999
$ws.Range("N141").Value = -35457.142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6096.8
$ws.Range("I32").Value = 5875.914
$ws.Range("J32").Value = 9031.429
$ws.Range("K32").Value = 5875.914
$ws.Range("L32").Value = 9031.429
$ws.Range("M32").Value = -5588.914
$ws.Range("N32").Value = -9605.429
$ws.Range("H45").Value = 10786304
$ws.Range("I45").Value = 13374601
$ws.Range("K45").Value = 13374601
$ws.Range("M45").Value = -13374224
$ws.Range("H122").Value = 8737.25
$ws.Range("I122").Value = 16510.5
$ws.Range("K122").Value = 49531.5
$ws.Range("M122").Value = -47081.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 38464244
$ws.Range("I105").Value = 2908.842
$ws.Range("K105").Value = 2908.842
$ws.Range("M105").Value = -1161.842

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2329.822
$ws.Range("I31").Value = 1761.2894
$ws.Range("J31").Value = 2947.0857
$ws.Range("K31").Value = 1761.2894
$ws.Range("L31").Value = 2947.0857
$ws.Range("M31").Value = -1466.2894
$ws.Range("N31").Value = -3537.0857
$ws.Range("H34").Value = 2329.822
$ws.Range("I34").Value = 1761.2894
$ws.Range("J34").Value = 2947.0857
$ws.Range("K34").Value = 1761.2894
$ws.Range("L34").Value = 2947.0857
$ws.Range("M34").Value = -1559.2894
$ws.Range("N34").Value = -3351.0857

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 5952.737
$ws.Range("I33").Value = 1152.6666
$ws.Range("J33").Value = 8168.154
$ws.Range("K33").Value = 6915.9996
$ws.Range("L33").Value = 49008.924
$ws.Range("M33").Value = -6632.9996
$ws.Range("N33").Value = -49574.924
$ws.Range("H35").Value = 999.8
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 999.8
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 2999.4
$ws.Range("N35").Value = -3575.4
$ws.Range("M35").ClearContents()
$ws.Range("H36").Value = 100
$ws.Range("I36").Value = 100
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 300
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -131
$ws.Range("N36").ClearContents()
$ws.Range("H69").Value = 11426.454
$ws.Range("I69").Value = 760.4
$ws.Range("J69").Value = 20314.834
$ws.Range("K69").Value = 2281.2
$ws.Range("L69").Value = 60944.50199999999
$ws.Range("M69").Value = -1470.2
$ws.Range("N69").Value = -62566.50199999999
$ws.Range("H72").Value = 11426.454
$ws.Range("I72").Value = 760.4
$ws.Range("J72").Value = 20314.834
$ws.Range("K72").Value = 6843.599999999999
$ws.Range("L72").Value = 182833.506
$ws.Range("M72").Value = -2787.599999999999
$ws.Range("N72").Value = -190945.506
$ws.Range("H124").Value = 2800
$ws.Range("I124").Value = 2000
$ws.Range("K124").Value = 6000
$ws.Range("M124").Value = -1090
$ws.Range("H131").Value = 700.25
$ws.Range("J131").Value = 746.8506
$ws.Range("L131").Value = 2240.5518
$ws.Range("N131").Value = -12320.5518

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4086.2727
$ws.Range("I70").Value = 4061.4736
$ws.Range("K70").Value = 4061.4736
$ws.Range("M70").Value = -3791.4736
$ws.Range("H73").Value = 4086.2727
$ws.Range("I73").Value = 4061.4736
$ws.Range("K73").Value = 4061.4736
$ws.Range("M73").Value = -3125.4736
$ws.Range("H97").Value = 1330
$ws.Range("I97").Value = 1255
$ws.Range("J97").Value = 1420
$ws.Range("K97").Value = 1255
$ws.Range("L97").Value = 1420
$ws.Range("M97").Value = -759
$ws.Range("N97").Value = -2412
$ws.Range("H102").Value = 894.5
$ws.Range("I102").Value = 782.2
$ws.Range("J102").Value = 1006.8
$ws.Range("K102").Value = 782.2
$ws.Range("L102").Value = 1006.8
$ws.Range("M102").Value = 839.8
$ws.Range("N102").Value = -4250.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 63026730
$ws.Range("I16").Value = 6494537
$ws.Range("J16").Value = 166669090
$ws.Range("K16").Value = 6494537
$ws.Range("L16").Value = 166669090
$ws.Range("M16").Value = -6494367
$ws.Range("N16").Value = -166669430
$ws.Range("H82").Value = 1240.1
$ws.Range("J82").Value = 1201
$ws.Range("L82").Value = 1201
$ws.Range("N82").Value = -1923
$ws.Range("H85").Value = 1240.1
$ws.Range("J85").Value = 1201
$ws.Range("L85").Value = 1201
$ws.Range("N85").Value = -3697
$ws.Range("H122").Value = 9302.775
$ws.Range("I122").Value = 11886.577
$ws.Range("J122").Value = 4504.2856
$ws.Range("K122").Value = 35659.731
$ws.Range("L122").Value = 13512.8568
$ws.Range("M122").Value = -33209.731
$ws.Range("N122").Value = -18412.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1710.4242
$ws.Range("I122").Value = 1055.5385
$ws.Range("J122").Value = 4142.857
$ws.Range("K122").Value = 3166.6155
$ws.Range("L122").Value = 12428.571
$ws.Range("M122").Value = -716.6155000000003
$ws.Range("N122").Value = -17328.571
$ws.Range("H123").Value = 19309.666
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 19309.666
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 19309.666
$ws.Range("N123").Value = -29109.666
$ws.Range("M123").ClearContents()
$ws.Range("H124").Value = 30000
$ws.Range("J124").Value = 30000
$ws.Range("L124").Value = 30000
$ws.Range("N124").Value = -39820
$ws.Range("H125").Value = 43072
$ws.Range("J125").Value = 43072
$ws.Range("L125").Value = 43072
$ws.Range("N125").Value = -52912
$ws.Range("H126").Value = 1154.0555
$ws.Range("I126").Value = 892.0625
$ws.Range("J126").Value = 3250
$ws.Range("K126").Value = 2676.1875
$ws.Range("L126").Value = 9750
$ws.Range("M126").Value = -206.1875
$ws.Range("N126").Value = -14690
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

